# Edit: Clean up Channel/Station labeling in Channel-width_LiDAR-GIS.xlsx
# - Column A: numeric channel id (e.g. 21) -> text "Channel <id>"
# - Column B: station code -> text "Station <code>" (with several data corrections:
#   zero-padded station numbers for Channel 19/21, and a handful of corrected
#   digits / BAF<->ORD labels that were mis-transcribed in the raw import)
# - Widen column A to fit the new "Channel NN" labels
# - Restore the active selection to F11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Cells.Item(2, 1).Value = "Channel 21"
$ws.Cells.Item(2, 2).Value = "Station 01 + 07.5 BAF"
$ws.Cells.Item(3, 1).Value = "Channel 21"
$ws.Cells.Item(3, 2).Value = "Station 01 + 49 ORD"
$ws.Cells.Item(4, 1).Value = "Channel 21"
$ws.Cells.Item(4, 2).Value = "Station 02 + 15 BAF"
$ws.Cells.Item(5, 1).Value = "Channel 21"
$ws.Cells.Item(5, 2).Value = "Station 03 + 80 BAF"
$ws.Cells.Item(6, 1).Value = "Channel 21"
$ws.Cells.Item(6, 2).Value = "Station 04 + 90 ORD"
$ws.Cells.Item(7, 1).Value = "Channel 21"
$ws.Cells.Item(7, 2).Value = "Station 06 + 36 BAF"
$ws.Cells.Item(8, 1).Value = "Channel 21"
$ws.Cells.Item(8, 2).Value = "Station 07 + 49.5 ORD"
$ws.Cells.Item(9, 1).Value = "Channel 21"
$ws.Cells.Item(9, 2).Value = "Station 07 + 66 BAF"
$ws.Cells.Item(10, 1).Value = "Channel 21"
$ws.Cells.Item(10, 2).Value = "Station 08 + 62 BAF"
$ws.Cells.Item(11, 1).Value = "Channel 21"
$ws.Cells.Item(11, 2).Value = "Station 09 + 99 ORD"
$ws.Cells.Item(12, 1).Value = "Channel 21"
$ws.Cells.Item(12, 2).Value = "Station 10 + 13 BAF"
$ws.Cells.Item(13, 1).Value = "Channel 21"
$ws.Cells.Item(13, 2).Value = "Station 11 + 12 ORD"
$ws.Cells.Item(14, 1).Value = "Channel 21"
$ws.Cells.Item(14, 2).Value = "Station 11 + 24 BAF"
$ws.Cells.Item(15, 1).Value = "Channel 21"
$ws.Cells.Item(15, 2).Value = "Station 11 + 63 ORD"
$ws.Cells.Item(16, 1).Value = "Channel 21"
$ws.Cells.Item(16, 2).Value = "Station 11 + 73 BAF"
$ws.Cells.Item(17, 1).Value = "Channel 19"
$ws.Cells.Item(17, 2).Value = "Station 01 + 04"
$ws.Cells.Item(18, 1).Value = "Channel 19"
$ws.Cells.Item(18, 2).Value = "Station 02 + 04"
$ws.Cells.Item(19, 1).Value = "Channel 19"
$ws.Cells.Item(19, 2).Value = "Station 02 + 52"
$ws.Cells.Item(20, 1).Value = "Channel 19"
$ws.Cells.Item(20, 2).Value = "Station 03 + 63"
$ws.Cells.Item(21, 1).Value = "Channel 19"
$ws.Cells.Item(21, 2).Value = "Station 03 + 98"
$ws.Cells.Item(22, 1).Value = "Channel 19"
$ws.Cells.Item(22, 2).Value = "Station 05 + 00"
$ws.Cells.Item(23, 1).Value = "Channel 19"
$ws.Cells.Item(23, 2).Value = "Station 05 + 09"
$ws.Cells.Item(24, 1).Value = "Channel 19"
$ws.Cells.Item(24, 2).Value = "Station 05 + 62"
$ws.Cells.Item(25, 1).Value = "Channel 19"
$ws.Cells.Item(25, 2).Value = "Station 06 + 13"
$ws.Cells.Item(26, 1).Value = "Channel 19"
$ws.Cells.Item(26, 2).Value = "Station 06 + 97"
$ws.Cells.Item(27, 1).Value = "Channel 19"
$ws.Cells.Item(27, 2).Value = "Station 07 + 30"
$ws.Cells.Item(28, 1).Value = "Channel 19"
$ws.Cells.Item(28, 2).Value = "Station 07 + 50"
$ws.Cells.Item(29, 1).Value = "Channel 19"
$ws.Cells.Item(29, 2).Value = "Station 07 + 97"
$ws.Cells.Item(30, 1).Value = "Channel 19"
$ws.Cells.Item(30, 2).Value = "Station 09 + 05"
$ws.Cells.Item(31, 1).Value = "Channel 19"
$ws.Cells.Item(31, 2).Value = "Station 09 + 26"
$ws.Cells.Item(32, 1).Value = "Channel 19"
$ws.Cells.Item(32, 2).Value = "Station 09 + 68"
$ws.Cells.Item(33, 1).Value = "Channel 19"
$ws.Cells.Item(33, 2).Value = "Station 10 + 98"
$ws.Cells.Item(34, 1).Value = "Channel 13"
$ws.Cells.Item(34, 2).Value = "Station 1 + 11.5 BAF"
$ws.Cells.Item(35, 1).Value = "Channel 13"
$ws.Cells.Item(35, 2).Value = "Station 1 + 38.5 ORD"
$ws.Cells.Item(36, 1).Value = "Channel 13"
$ws.Cells.Item(36, 2).Value = "Station 1 + 71 BAF"
$ws.Cells.Item(37, 1).Value = "Channel 13"
$ws.Cells.Item(37, 2).Value = "Station 1 + 77 ORD"
$ws.Cells.Item(38, 1).Value = "Channel 13"
$ws.Cells.Item(38, 2).Value = "Station 1 + 99 ORD"
$ws.Cells.Item(39, 1).Value = "Channel 13"
$ws.Cells.Item(39, 2).Value = "Station 2 + 16 BAF"
$ws.Cells.Item(40, 1).Value = "Channel 13"
$ws.Cells.Item(40, 2).Value = "Station 2 + 57 ORD"
$ws.Cells.Item(41, 1).Value = "Channel 13"
$ws.Cells.Item(41, 2).Value = "Station 2 + 67 BAF"
$ws.Cells.Item(42, 1).Value = "Channel 13"
$ws.Cells.Item(42, 2).Value = "Station 3 + 09 ORD"
$ws.Cells.Item(43, 1).Value = "Channel 13"
$ws.Cells.Item(43, 2).Value = "Station 3 + 31 BAF"
$ws.Cells.Item(44, 1).Value = "Channel 13"
$ws.Cells.Item(44, 2).Value = "Station 3 + 65 ORD"
$ws.Cells.Item(45, 1).Value = "Channel 13"
$ws.Cells.Item(45, 2).Value = "Station 3 + 95.5 ORD"
$ws.Cells.Item(46, 1).Value = "Channel 13"
$ws.Cells.Item(46, 2).Value = "Station 4 + 14.5 BAF"
$ws.Cells.Item(47, 1).Value = "Channel 13"
$ws.Cells.Item(47, 2).Value = "Station 4 + 36 ORD"
$ws.Cells.Item(48, 1).Value = "Channel 13"
$ws.Cells.Item(48, 2).Value = "Station 4 + 50 BAF"
$ws.Cells.Item(49, 1).Value = "Channel 13"
$ws.Cells.Item(49, 2).Value = "Station 4 + 78 ORD"
$ws.Cells.Item(50, 1).Value = "Channel 12"
$ws.Cells.Item(50, 2).Value = "Station 1 + 04"
$ws.Cells.Item(51, 1).Value = "Channel 12"
$ws.Cells.Item(51, 2).Value = "Station 1 + 52"
$ws.Cells.Item(52, 1).Value = "Channel 12"
$ws.Cells.Item(52, 2).Value = "Station 1 + 83"
$ws.Cells.Item(53, 1).Value = "Channel 12"
$ws.Cells.Item(53, 2).Value = "Station 2 + 63"
$ws.Cells.Item(54, 1).Value = "Channel 12"
$ws.Cells.Item(54, 2).Value = "Station 2 + 98"
$ws.Cells.Item(55, 1).Value = "Channel 12"
$ws.Cells.Item(55, 2).Value = "Station 3 + 93"
$ws.Cells.Item(56, 1).Value = "Channel 12"
$ws.Cells.Item(56, 2).Value = "Station 4 + 00"
$ws.Cells.Item(57, 1).Value = "Channel 12"
$ws.Cells.Item(57, 2).Value = "Station 4 + 09"
$ws.Cells.Item(58, 1).Value = "Channel 12"
$ws.Cells.Item(58, 2).Value = "Station 4 + 56"
$ws.Cells.Item(59, 1).Value = "Channel 12"
$ws.Cells.Item(59, 2).Value = "Station 4 + 65"
$ws.Cells.Item(60, 1).Value = "Channel 12"
$ws.Cells.Item(60, 2).Value = "Station 5 + 13"
$ws.Cells.Item(61, 1).Value = "Channel 12"
$ws.Cells.Item(61, 2).Value = "Station 5 + 48"
$ws.Cells.Item(62, 1).Value = "Channel 12"
$ws.Cells.Item(62, 2).Value = "Station 6 + 30"
$ws.Cells.Item(63, 1).Value = "Channel 12"
$ws.Cells.Item(63, 2).Value = "Station 6 + 93"


# Column A needs to be a bit wider to comfortably fit "Channel NN" labels
$ws.Columns.Item(1).ColumnWidth = 10.83

# Restore selection to F11 (bottom-left frozen pane)
$ws.Range("F11").Select()
